# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" sheet with newer totals for a handful of
# countries and re-sort the affected rows by total cases (column B,
# descending) so the table stays ordered, matching how the source feed
# re-ran after the 08:59 -> 09:29 refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 09:29"

# --- Simple value refreshes (no re-sort needed) ------------------------------

# Row 8: Alemania
$ws.Range("B8").Value = 53340
$ws.Range("C8").Value = 2469
$ws.Range("D8").Value = 6658
$ws.Range("E8").Value = 46287
$ws.Range("F8").Value = 1581
$ws.Range("G8").Value = 44
$ws.Range("H8").Value = 395

# Row 15: Austria
$ws.Range("B15").Value = 7712
$ws.Range("C15").Value = 15
$ws.Range("D15").Value = 225
$ws.Range("E15").Value = 7429
$ws.Range("F15").Value = 128
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 58

# Row 61: Barein
$ws.Range("B61").Value = 473
$ws.Range("C61").Value = 7
$ws.Range("D61").Value = 254
$ws.Range("E61").Value = 215
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 4

# --- Rows whose countries swap order after re-sorting by total cases --------

# Polonia (1436) now outranks Pakistan (1408); row 33 becomes Polonia,
# row 34 becomes Pakistan.
$ws.Range("A33").Value = "Polonia"
$ws.Range("B33").Value = 1436
$ws.Range("C33").Value = 47
$ws.Range("D33").Value = 7
$ws.Range("E33").Value = 1413
$ws.Range("F33").Value = 3
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 16

$ws.Range("A34").Value = "Pakistan"
$ws.Range("B34").Value = 1408
$ws.Range("C34").Value = 35
$ws.Range("D34").Value = 25
$ws.Range("E34").Value = 1372
$ws.Range("F34").Value = 7
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 11

# Croacia (635) now outranks Eslovenia (632); row 52 becomes Croacia,
# row 53 becomes Eslovenia.
$ws.Range("A52").Value = "Croacia"
$ws.Range("B52").Value = 635
$ws.Range("C52").Value = 49
$ws.Range("D52").Value = 45
$ws.Range("E52").Value = 586
$ws.Range("F52").Value = 14
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 4

$ws.Range("A53").Value = "Eslovenia"
$ws.Range("B53").Value = 632
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 10
$ws.Range("E53").Value = 613
$ws.Range("F53").Value = 14
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 9

# Nigeria (81) jumps above Cuba/Bolivia/Guadalupe; row 113 becomes
# Nigeria, and Cuba/Bolivia/Guadalupe each shift down one row.
$ws.Range("A113").Value = "Nigeria"
$ws.Range("B113").Value = 81
$ws.Range("C113").Value = 11
$ws.Range("D113").Value = 3
$ws.Range("E113").Value = 77
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 1

$ws.Range("A114").Value = "Cuba"
$ws.Range("B114").Value = 80
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 4
$ws.Range("E114").Value = 74
$ws.Range("F114").Value = 2
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 2

$ws.Range("A115").Value = "Bolivia"
$ws.Range("B115").Value = 74
$ws.Range("C115").Value = 13
$ws.Range("D115").Value = 0
$ws.Range("E115").Value = 74
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 0

$ws.Range("A116").Value = "Guadalupe"
$ws.Range("B116").Value = 73
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 0
$ws.Range("E116").Value = 72
$ws.Range("F116").Value = 4
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 1
